$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3127293
$ws.Range("I40").Value = 6946000
$ws.Range("J40").Value = 2896.3635
$ws.Range("K40").Value = 6946000
$ws.Range("L40").Value = 2896.3635
$ws.Range("M40").Value = -6945825
$ws.Range("N40").Value = -3246.3635
$ws.Range("H64").Value = 3085.5144
$ws.Range("I64").Value = 2817.8635
$ws.Range("J64").Value = 3538.4614
$ws.Range("K64").Value = 2817.8635
$ws.Range("L64").Value = 3538.4614
$ws.Range("M64").Value = -2569.8635
$ws.Range("N64").Value = -4034.4614
$ws.Range("H67").Value = 3085.5144
$ws.Range("I67").Value = 2817.8635
$ws.Range("J67").Value = 3538.4614
$ws.Range("K67").Value = 2817.8635
$ws.Range("L67").Value = 3538.4614
$ws.Range("M67").Value = -1959.8635
$ws.Range("N67").Value = -5254.4614
$ws.Range("H76").Value = 15159340
$ws.Range("I76").Value = 9289.529
$ws.Range("J76").Value = 66669510
$ws.Range("K76").Value = 9289.529
$ws.Range("L76").Value = 66669510
$ws.Range("M76").Value = -8974.529
$ws.Range("N76").Value = -66670140
$ws.Range("H79").Value = 15159340
$ws.Range("I79").Value = 9289.529
$ws.Range("J79").Value = 66669510
$ws.Range("K79").Value = 9289.529
$ws.Range("L79").Value = 66669510
$ws.Range("M79").Value = -8197.529
$ws.Range("N79").Value = -66671694
$ws.Range("H132").Value = 7250989.5
$ws.Range("I132").Value = 8337424
$ws.Range("J132").Value = 8093.778
$ws.Range("K132").Value = 25012272
$ws.Range("L132").Value = 24281.334
$ws.Range("M132").Value = -25009742
$ws.Range("N132").Value = -29341.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 48706.332
$ws.Range("I45").Value = 84161.086
$ws.Range("J45").Value = 1433.3334
$ws.Range("K45").Value = 84161.086
$ws.Range("L45").Value = 1433.3334
$ws.Range("M45").Value = -83784.086
$ws.Range("N45").Value = -2187.3334
$ws.Range("H63").Value = 2177.5557
$ws.Range("I63").Value = 2099.875
$ws.Range("K63").Value = 2099.875
$ws.Range("M63").Value = -1413.875
$ws.Range("H66").Value = 2177.5557
$ws.Range("I66").Value = 2099.875
$ws.Range("K66").Value = 10499.375
$ws.Range("M66").Value = -7067.375
$ws.Range("H88").Value = 2515.4583
$ws.Range("I88").Value = 3352.4443
$ws.Range("J88").Value = 2013.2667
$ws.Range("K88").Value = 3352.4443
$ws.Range("L88").Value = 2013.2667
$ws.Range("M88").Value = -2946.4443
$ws.Range("N88").Value = -2825.2667
$ws.Range("H91").Value = 2515.4583
$ws.Range("I91").Value = 3352.4443
$ws.Range("J91").Value = 2013.2667
$ws.Range("K91").Value = 3352.4443
$ws.Range("L91").Value = 2013.2667
$ws.Range("M91").Value = -1948.4443
$ws.Range("N91").Value = -4821.2667
$ws.Range("H122").Value = 1526.0526
$ws.Range("I122").Value = 1499.7222
$ws.Range("K122").Value = 4499.1666
$ws.Range("M122").Value = -2049.1666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1013411.06
$ws.Range("I86").Value = 2772.4614
$ws.Range("J86").Value = 2327241.2
$ws.Range("K86").Value = 2772.4614
$ws.Range("L86").Value = 2327241.2
$ws.Range("M86").Value = -1649.4614
$ws.Range("N86").Value = -2329487.2
$ws.Range("H89").Value = 1013411.06
$ws.Range("I89").Value = 2772.4614
$ws.Range("J89").Value = 2327241.2
$ws.Range("K89").Value = 13862.307
$ws.Range("L89").Value = 11636206
$ws.Range("M89").Value = -8246.307000000001
$ws.Range("N89").Value = -11647438
$ws.Range("H105").Value = 37039330
$ws.Range("I105").Value = 2369.95
$ws.Range("K105").Value = 2369.95
$ws.Range("M105").Value = -622.9499999999998
$ws.Range("H134").Value = 1739055.4
$ws.Range("I134").Value = 1012.6667
$ws.Range("K134").Value = 3038.0001
$ws.Range("M134").Value = -503.0001000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1408.439
$ws.Range("I16").Value = 739.04346
$ws.Range("J16").Value = 2263.7778
$ws.Range("K16").Value = 739.04346
$ws.Range("L16").Value = 2263.7778
$ws.Range("M16").Value = -452.04346
$ws.Range("N16").Value = -2837.7778
$ws.Range("H31").Value = 3729.0366
$ws.Range("I31").Value = 2713.8948
$ws.Range("J31").Value = 6043.56
$ws.Range("K31").Value = 2713.8948
$ws.Range("L31").Value = 6043.56
$ws.Range("M31").Value = -2418.8948
$ws.Range("N31").Value = -6633.56
$ws.Range("H34").Value = 3729.0366
$ws.Range("I34").Value = 2713.8948
$ws.Range("J34").Value = 6043.56
$ws.Range("K34").Value = 2713.8948
$ws.Range("L34").Value = 6043.56
$ws.Range("M34").Value = -2511.8948
$ws.Range("N34").Value = -6447.56
$ws.Range("H62").Value = 4307.1816
$ws.Range("I62").Value = 2949.5
$ws.Range("K62").Value = 2949.5
$ws.Range("M62").Value = -2325.5
$ws.Range("H65").Value = 4307.1816
$ws.Range("I65").Value = 2949.5
$ws.Range("K65").Value = 14747.5
$ws.Range("M65").Value = -11627.5
$ws.Range("H104").Value = 36755
$ws.Range("I104").Value = 30000
$ws.Range("J104").Value = 40132.5
$ws.Range("K104").Value = 30000
$ws.Range("L104").Value = 40132.5
$ws.Range("M104").Value = -27379
$ws.Range("N104").Value = -45374.5
$ws.Range("H113").Value = 1408.439
$ws.Range("I113").Value = 739.04346
$ws.Range("J113").Value = 2263.7778
$ws.Range("K113").Value = 739.04346
$ws.Range("L113").Value = 2263.7778
$ws.Range("M113").Value = 1430.95654
$ws.Range("N113").Value = -6603.7778
$ws.Range("H132").Value = 16673522
$ws.Range("I132").Value = 1365.2858
$ws.Range("J132").Value = 55575220
$ws.Range("K132").Value = 4095.8574
$ws.Range("L132").Value = 166725660
$ws.Range("M132").Value = -1565.8574
$ws.Range("N132").Value = -166730720

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 34413.8
$ws.Range("I103").Value = 475
$ws.Range("J103").Value = 46755.184
$ws.Range("K103").Value = 1425
$ws.Range("L103").Value = 140265.552
$ws.Range("M103").Value = -546
$ws.Range("N103").Value = -142023.552
$ws.Range("H113").Value = 12222684
$ws.Range("I113").Value = 10417141
$ws.Range("J113").Value = 14286162
$ws.Range("K113").Value = 31251423
$ws.Range("L113").Value = 42858486
$ws.Range("M113").Value = -31249253
$ws.Range("N113").Value = -42862826
$ws.Range("H131").Value = 699.35
$ws.Range("J131").Value = 778.3125
$ws.Range("L131").Value = 2334.9375
$ws.Range("N131").Value = -12414.9375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4240.7085
$ws.Range("I70").Value = 4225.1055
$ws.Range("J70").Value = 4300
$ws.Range("K70").Value = 4225.1055
$ws.Range("L70").Value = 4300
$ws.Range("M70").Value = -3955.1055
$ws.Range("N70").Value = -4840
$ws.Range("H73").Value = 4240.7085
$ws.Range("I73").Value = 4225.1055
$ws.Range("J73").Value = 4300
$ws.Range("K73").Value = 4225.1055
$ws.Range("L73").Value = 4300
$ws.Range("M73").Value = -3289.1055
$ws.Range("N73").Value = -6172
$ws.Range("H80").Value = 12502175
$ws.Range("I80").Value = 2483.8333
$ws.Range("J80").Value = 50001250
$ws.Range("K80").Value = 2483.8333
$ws.Range("L80").Value = 50001250
$ws.Range("M80").Value = -1485.8333
$ws.Range("N80").Value = -50003246
$ws.Range("H83").Value = 12502175
$ws.Range("I83").Value = 2483.8333
$ws.Range("J83").Value = 50001250
$ws.Range("K83").Value = 12419.1665
$ws.Range("L83").Value = 250006250
$ws.Range("M83").Value = -7427.166499999999
$ws.Range("N83").Value = -250016234
$ws.Range("H132").Value = 10991.091
$ws.Range("I132").Value = 872.5
$ws.Range("J132").Value = 16773.143
$ws.Range("K132").Value = 2617.5
$ws.Range("L132").Value = 50319.429
$ws.Range("M132").Value = -87.5
$ws.Range("N132").Value = -55379.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1520.5264
$ws.Range("H71").Value = 1520.5264
$ws.Range("H132").Value = 16671142
$ws.Range("I132").Value = 31251848
$ws.Range("J132").Value = 7478.6787
$ws.Range("K132").Value = 93755544
$ws.Range("L132").Value = 22436.0361
$ws.Range("M132").Value = -93753014
$ws.Range("N132").Value = -27496.0361

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 18042.762
$ws.Range("I132").Value = 23937.045
$ws.Range("J132").Value = 5986.273
$ws.Range("K132").Value = 71811.13499999999
$ws.Range("L132").Value = 17958.819
$ws.Range("M132").Value = -69281.13499999999
$ws.Range("N132").Value = -23018.819
$ws.Range("H136").Value = 4119.525
$ws.Range("I136").Value = 7980.067
$ws.Range("K136").Value = 23940.201
$ws.Range("M136").Value = -21390.201
